# Add season record columns (Wins, Losses, Ties) to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row/column of the existing table so we append
# the new "Wins"/"Losses"/"Ties" columns right after it (column AC -> AD:AF).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count
$newColStart = $lastCol + 1

# Header row: copy style of the existing header cell (A1) so the new
# headers match the look of the rest of the header row (bold, bordered,
# centered).
$headerSrc = $ws.Cells.Item(1, 1)

$winsHeader = $ws.Cells.Item(1, $newColStart)
$winsHeader.Value = "Wins"
$headerSrc.Copy()
$winsHeader.PasteSpecial(-4122)
$winsHeader.Value = "Wins"

$lossesHeader = $ws.Cells.Item(1, $newColStart + 1)
$lossesHeader.Value = "Losses"
$headerSrc.Copy()
$lossesHeader.PasteSpecial(-4122)
$lossesHeader.Value = "Losses"

$tiesHeader = $ws.Cells.Item(1, $newColStart + 2)
$tiesHeader.Value = "Ties"
$headerSrc.Copy()
$tiesHeader.PasteSpecial(-4122)
$tiesHeader.Value = "Ties"

$excel.CutCopyMode = 0

# Season record values, same for every player row in this sheet (2001
# Cincinnati Reds: 66 wins, 96 losses, 0 ties).
$wins = 66
$losses = 96
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newColStart).Value = $wins
    $ws.Cells.Item($r, $newColStart + 1).Value = $losses
    $ws.Cells.Item($r, $newColStart + 2).Value = $ties
}
